# Apply the "Atas - 2020 a 2023 - Dados" update:
#  - fill in Km de linha / Capacidade / Instalação / GD data for meetings 260-275
#  - add meeting 275 (new row 51), shifting the footnote row down
#  - mark meeting 262 as extraordinary ("262*") since it has no data
#  - apply thousands-separator formatting (#,##0) across the data columns C:F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 51 so the blank spacer rows + footnote shift down
#    (old row 53 footnote becomes row 54), matching the new meeting (275) row.
$ws.Rows("51:51").Insert()

# 2. Apply "#,##0" number format across the whole data block C2:F50 first -
#    this both updates the existing cells' styles and creates the still-empty
#    cells (e.g. F2:F4, C36:F50) with the correctly-formatted blank style.
$ws.Range("C2:F50").NumberFormat = "#,##0"

# 3. Fill in the data that was missing for meetings 260-274 (rows 36-50).
$ws.Range("C36").Value = 1122
$ws.Range("D36").Value = 1213
$ws.Range("E36").Value = 450
$ws.Range("F36").Value = "-"

$ws.Range("C37").Value = 482
$ws.Range("D37").Value = 687
$ws.Range("E37").Value = "-"
$ws.Range("F37").Value = "-"

# Meeting 262 never happened / has no minutes -> mark as extraordinary like
# the other "*" meetings, and blank out its data row with "-".
$ws.Range("B38").Value = "262*"
$ws.Range("C38").Value = "-"
$ws.Range("D38").Value = "-"
$ws.Range("E38").Value = "-"
$ws.Range("F38").Value = "-"

$ws.Range("C39").Value = 516
$ws.Range("D39").Value = 566
$ws.Range("E39").Value = 2396
$ws.Range("F39").Value = "-"

$ws.Range("C40").Value = 347
$ws.Range("D40").Value = 409
$ws.Range("E40").Value = 2226
$ws.Range("F40").Value = "-"

$ws.Range("C41").Value = 223
$ws.Range("D41").Value = 645
$ws.Range("E41").Value = 400
$ws.Range("F41").Value = "-"

$ws.Range("C42").Value = 596
$ws.Range("D42").Value = 1012
$ws.Range("E42").Value = 1250
$ws.Range("F42").Value = "-"

$ws.Range("C43").Value = 158
$ws.Range("D43").Value = 303
$ws.Range("E43").Value = 700
$ws.Range("F43").Value = "-"

$ws.Range("C44").Value = 510
$ws.Range("D44").Value = 544
$ws.Range("E44").Value = "-"
$ws.Range("F44").Value = "-"

$ws.Range("C45").Value = 650
$ws.Range("D45").Value = 139
$ws.Range("E45").Value = "-"
$ws.Range("F45").Value = "-"

$ws.Range("C46").Value = 1387
$ws.Range("D46").Value = 244
$ws.Range("E46").Value = "-"
$ws.Range("F46").Value = "-"

$ws.Range("C47").Value = 941
$ws.Range("D47").Value = 152
$ws.Range("E47").Value = 300
$ws.Range("F47").Value = "-"

$ws.Range("C48").Value = 1030
$ws.Range("D48").Value = 1223
$ws.Range("E48").Value = 750
$ws.Range("F48").Value = "-"

$ws.Range("C49").Value = "-"
$ws.Range("D49").Value = "-"
$ws.Range("E49").Value = "-"
$ws.Range("F49").Value = "-"

$ws.Range("C50").Value = 2138
$ws.Range("D50").Value = 455
$ws.Range("E50").Value = 1866
$ws.Range("F50").Value = 864

# 4. Row 35 (meeting 259*) - mark the whole data row with "-" placeholders.
$ws.Range("C35").Value = "-"
$ws.Range("D35").Value = "-"
$ws.Range("E35").Value = "-"
$ws.Range("F35").Value = "-"

# 5. Row 50 is no longer the final row of the table, so it loses the
#    bottom-border "closing" style and becomes a normal data row.
$ws.Range("A50").Style = $ws.Range("A49").Style
$ws.Range("B50").Style = $ws.Range("B49").Style
$ws.Range("C50").Style = $ws.Range("C49").Style
$ws.Range("D50").Style = $ws.Range("D49").Style
$ws.Range("E50").Style = $ws.Range("E49").Style
$ws.Range("F50").Style = $ws.Range("F49").Style
$ws.Range("G50").Style = $ws.Range("G49").Style

# 6. New row 51: meeting 275, closes out the table so it gets the
#    bottom-border "closing" style that row 50 used to have.
$ws.Range("A51").Value = 2023
$ws.Range("B51").Value = 275
$ws.Range("C51").Value = 749
$ws.Range("D51").Value = 1096
$ws.Range("E51").Value = "-"

$ws.Range("A51").Style = $ws.Range("A48").Style
$ws.Range("B51").Style = $ws.Range("B48").Style
$ws.Range("C51").Style = $ws.Range("C48").Style
$ws.Range("D51").Style = $ws.Range("D48").Style
$ws.Range("E51").Style = $ws.Range("E48").Style
$ws.Range("F51").Style = $ws.Range("F48").Style
$ws.Range("G51").Style = $ws.Range("G48").Style

# 7. Tidy up the view: land on the newly-added row like the author did.
$ws.Range("A51").Select()

Write-Host "Edit applied"
